$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 8) representing a "None" buff - grants the team nothing.
$ws.Range("A8").Value = "None"
$ws.Range("B8").Value = "No Title"
$ws.Range("C8").Value = "grants the team nothing"

# Numeric stat columns D:O all default to 0 for this buff.
$ws.Range("D8:O8").Value = 0

# Update view state to reflect how the workbook was left after editing.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D13").Select()
